$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 3) to the report sheet, mirroring the structure of row 2.
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "حسن "
$ws.Range("C3").Value = "'2222"
$ws.Range("D3").Value = "ايتا"
$ws.Range("E3").Value = "الرحلة 2"
$ws.Range("F3").Value = "C3"
$ws.Range("G3").Value = "NRC"
$ws.Range("H3").Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠١:٥٠:٥٣ م"
